$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2921.3125
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2921.3125
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8763.9375
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = -9099.9375

$ws.Range("H32").Value = 6499.6665
$ws.Range("I32").Value = 4999.6665
$ws.Range("K32").Value = 4999.6665
$ws.Range("M32").Value = -4673.6665

$ws.Range("H94").Value = 7359.2
$ws.Range("I94").Value = 7359.2
$ws.Range("K94").Value = 7359.2
$ws.Range("M94").Value = -6908.2

$ws.Range("H98").Value = 3641.7334
$ws.Range("I98").Value = 3782.7
$ws.Range("J98").Value = 3359.8
$ws.Range("K98").Value = 3782.7
$ws.Range("L98").Value = 3359.8
$ws.Range("M98").Value = -2284.7
$ws.Range("N98").Value = -6355.8

$ws.Range("H122").Value = 3641.7334
$ws.Range("I122").Value = 3782.7
$ws.Range("J122").Value = 3359.8
$ws.Range("K122").Value = 11348.1
$ws.Range("L122").Value = 10079.4
$ws.Range("M122").Value = -8898.099999999999
$ws.Range("N122").Value = -14979.4

$ws.Range("H132").Value = 2372.4443
$ws.Range("I132").Value = 2290.6
$ws.Range("J132").Value = 2474.75
$ws.Range("K132").Value = 6871.799999999999
$ws.Range("L132").Value = 7424.25
$ws.Range("M132").Value = -4341.799999999999
$ws.Range("N132").Value = -12484.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3326.5557
$ws.Range("I74").Value = 2264.4546
$ws.Range("J74").Value = 4995.5713
$ws.Range("K74").Value = 2264.4546
$ws.Range("L74").Value = 4995.5713
$ws.Range("M74").Value = -1390.4546
$ws.Range("N74").Value = -6743.5713

$ws.Range("H77").Value = 3326.5557
$ws.Range("I77").Value = 2264.4546
$ws.Range("J77").Value = 4995.5713
$ws.Range("K77").Value = 11322.273
$ws.Range("L77").Value = 24977.8565
$ws.Range("M77").Value = -6954.273000000001
$ws.Range("N77").Value = -33713.85649999999

$ws.Range("H122").Value = 5360.702
$ws.Range("I122").Value = 5310.65
$ws.Range("J122").Value = 5646.7144
$ws.Range("K122").Value = 15931.95
$ws.Range("L122").Value = 16940.1432
$ws.Range("M122").Value = -13481.95
$ws.Range("N122").Value = -21840.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3971.3572
$ws.Range("I94").Value = 2057.5
$ws.Range("K94").Value = 2057.5
$ws.Range("M94").Value = -1606.5

$ws.Range("H105").Value = 1859.4736
$ws.Range("I105").Value = 1851.6666
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1851.6666
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -104.6666
$ws.Range("N105").Value = -5494

$ws.Range("H134").Value = 2476.4092
$ws.Range("I134").Value = 1975.2858
$ws.Range("K134").Value = 5925.857400000001
$ws.Range("M134").Value = -3390.857400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1200.3334
$ws.Range("J16").Value = 1198.5
$ws.Range("L16").Value = 1198.5
$ws.Range("N16").Value = -1772.5

$ws.Range("H31").Value = 1135.5555
$ws.Range("I31").Value = 1135.5555
$ws.Range("K31").Value = 1135.5555
$ws.Range("M31").Value = -840.5554999999999

$ws.Range("H34").Value = 1135.5555
$ws.Range("I34").Value = 1135.5555
$ws.Range("K34").Value = 1135.5555
$ws.Range("M34").Value = -933.5554999999999

$ws.Range("H62").Value = 6324
$ws.Range("I62").Value = 6324
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 6324
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -5700
$ws.Range("N62").Value = $null

$ws.Range("H65").Value = 6324
$ws.Range("I65").Value = 6324
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 31620
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -28500
$ws.Range("N65").Value = $null

$ws.Range("H107").Value = 889.29266
$ws.Range("I107").Value = 499.0625
$ws.Range("J107").Value = 1139.04
$ws.Range("K107").Value = 499.0625
$ws.Range("L107").Value = 1139.04
$ws.Range("M107").Value = 1420.9375
$ws.Range("N107").Value = -4979.04

$ws.Range("H113").Value = 1200.3334
$ws.Range("J113").Value = 1198.5
$ws.Range("L113").Value = 1198.5
$ws.Range("N113").Value = -5538.5

$ws.Range("H132").Value = 1880.24
$ws.Range("I132").Value = 1668.5
$ws.Range("K132").Value = 5005.5
$ws.Range("M132").Value = -2475.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 3657.8
$ws.Range("J80").Value = 3097.3333
$ws.Range("L80").Value = 9291.999899999999
$ws.Range("N80").Value = -11163.9999

$ws.Range("H83").Value = 3657.8
$ws.Range("J83").Value = 3097.3333
$ws.Range("L83").Value = 27875.9997
$ws.Range("N83").Value = -37235.9997

$ws.Range("H114").Value = 23455.1
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 23455.1
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 70365.29999999999
$ws.Range("M114").Value = $null
$ws.Range("N114").Value = -76873.29999999999

$ws.Range("H129").Value = 4143.5
$ws.Range("J129").Value = 2836.9
$ws.Range("L129").Value = 8510.700000000001
$ws.Range("N129").Value = -18510.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3188.5715
$ws.Range("I102").Value = 3217.6667
$ws.Range("J102").Value = 3014
$ws.Range("K102").Value = 3217.6667
$ws.Range("L102").Value = 3014
$ws.Range("M102").Value = -1595.6667
$ws.Range("N102").Value = -6258

$ws.Range("H113").Value = 2578.625
$ws.Range("I113").Value = 2428
$ws.Range("K113").Value = 2428
$ws.Range("M113").Value = -258

$ws.Range("H122").Value = 1319
$ws.Range("I122").Value = 1326
$ws.Range("K122").Value = 3978
$ws.Range("M122").Value = -1528

$ws.Range("H132").Value = 2890.8823
$ws.Range("I132").Value = 2587.5386
$ws.Range("J132").Value = 3876.75
$ws.Range("K132").Value = 7762.6158
$ws.Range("L132").Value = 11630.25
$ws.Range("M132").Value = -5232.6158
$ws.Range("N132").Value = -16690.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2650.4443
$ws.Range("I93").Value = 1669.8334
$ws.Range("K93").Value = 1669.8334
$ws.Range("M93").Value = -421.8334

$ws.Range("H132").Value = 4979.4644
$ws.Range("I132").Value = 4022.6667
$ws.Range("J132").Value = 7849.857
$ws.Range("K132").Value = 12068.0001
$ws.Range("L132").Value = 23549.571
$ws.Range("M132").Value = -9538.000100000001
$ws.Range("N132").Value = -28609.571

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 21288.375
$ws.Range("J54").Value = 21000
$ws.Range("L54").Value = 21000
$ws.Range("N54").Value = -22040

Write-Host "Applied Zodiark_Profits value updates"
